$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4625
$ws.Range("I62").Value = 4333.3335
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 4333.3335
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -3709.3335
$ws.Range("N62").Value = -6748

$ws.Range("H65").Value = 4625
$ws.Range("I65").Value = 4333.3335
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 21666.6675
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -18546.6675
$ws.Range("N65").Value = -33740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1209.52
$ws.Range("I2").Value = 1072.8235
$ws.Range("K2").Value = 1072.8235
$ws.Range("M2").Value = -959.8235

$ws.Range("H32").Value = 13128.23
$ws.Range("I32").Value = 6043.864
$ws.Range("J32").Value = 26880.234
$ws.Range("K32").Value = 6043.864
$ws.Range("L32").Value = 26880.234
$ws.Range("M32").Value = -5756.864
$ws.Range("N32").Value = -27454.234

$ws.Range("H45").Value = 1513.6
$ws.Range("I45").Value = 1655.8889
$ws.Range("J45").Value = 1147.7142
$ws.Range("K45").Value = 1655.8889
$ws.Range("L45").Value = 1147.7142
$ws.Range("M45").Value = -1278.8889
$ws.Range("N45").Value = -1901.7142

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H110").Value = 2012.8182
$ws.Range("I110").Value = 2067.625
$ws.Range("K110").Value = 2067.625
$ws.Range("M110").Value = -22.625

$ws.Range("H111").Value = 32000
$ws.Range("J111").Value = 32000
$ws.Range("L111").Value = 32000
$ws.Range("N111").Value = -40180

$ws.Range("H116").Value = 1209.52
$ws.Range("I116").Value = 1072.8235
$ws.Range("K116").Value = 1072.8235
$ws.Range("M116").Value = 1221.1765

$ws.Range("H122").Value = 3893.3333
$ws.Range("I122").Value = 3971.2
$ws.Range("J122").Value = 3737.6
$ws.Range("K122").Value = 11913.6
$ws.Range("L122").Value = 11212.8
$ws.Range("M122").Value = -9463.599999999999
$ws.Range("N122").Value = -16112.8

$ws.Range("H125").Value = 44003.75
$ws.Range("J125").Value = 44003.75
$ws.Range("L125").Value = 44003.75
$ws.Range("N125").Value = -53843.75

$ws.Range("H132").Value = 23095.043
$ws.Range("I132").Value = 25743.219
$ws.Range("J132").Value = 4999.1665
$ws.Range("K132").Value = 77229.65700000001
$ws.Range("L132").Value = 14997.4995
$ws.Range("M132").Value = -74699.65700000001
$ws.Range("N132").Value = -20057.4995

$ws.Range("H133").Value = 23128.834
$ws.Range("J133").Value = 23128.834
$ws.Range("L133").Value = 23128.834
$ws.Range("N133").Value = -28188.834

$ws.Range("H134").Value = 25718.285
$ws.Range("J134").Value = 25718.285
$ws.Range("L134").Value = 25718.285
$ws.Range("N134").Value = -35858.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1209.52
$ws.Range("I3").Value = 1072.8235
$ws.Range("K3").Value = 1072.8235
$ws.Range("M3").Value = -958.8235

$ws.Range("H99").Value = 1757.32
$ws.Range("I99").Value = 1464.2142
$ws.Range("J99").Value = 2130.3635
$ws.Range("K99").Value = 1464.2142
$ws.Range("L99").Value = 2130.3635
$ws.Range("M99").Value = 33.78580000000011
$ws.Range("N99").Value = -5126.363499999999

$ws.Range("H134").Value = 669086.5600000001
$ws.Range("I134").Value = 1250931
$ws.Range("J134").Value = 4121.4287
$ws.Range("K134").Value = 3752793
$ws.Range("L134").Value = 12364.2861
$ws.Range("M134").Value = -3750258
$ws.Range("N134").Value = -17434.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2154.1296
$ws.Range("I31").Value = 1382
$ws.Range("J31").Value = 2684.9688
$ws.Range("K31").Value = 1382
$ws.Range("L31").Value = 2684.9688
$ws.Range("M31").Value = -1087
$ws.Range("N31").Value = -3274.9688

$ws.Range("H34").Value = 2154.1296
$ws.Range("I34").Value = 1382
$ws.Range("J34").Value = 2684.9688
$ws.Range("K34").Value = 1382
$ws.Range("L34").Value = 2684.9688
$ws.Range("M34").Value = -1180
$ws.Range("N34").Value = -3088.9688

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 749.9
$ws.Range("I97").Value = 755.44446
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 2266.33338
$ws.Range("L97").Value = 2100
$ws.Range("M97").Value = -1770.33338
$ws.Range("N97").Value = -3092

$ws.Range("H121").Value = 5485.4346
$ws.Range("I121").Value = 6208.647
$ws.Range("J121").Value = 5061.483
$ws.Range("K121").Value = 18625.941
$ws.Range("L121").Value = 15184.449
$ws.Range("M121").Value = -17315.941
$ws.Range("N121").Value = -17804.449

$ws.Range("H131").Value = 875.9091
$ws.Range("I131").Value = 415
$ws.Range("J131").Value = 1048.75
$ws.Range("K131").Value = 1245
$ws.Range("L131").Value = 3146.25
$ws.Range("M131").Value = 3795
$ws.Range("N131").Value = -13226.25

$ws.Range("H141").Value = 1213.7222
$ws.Range("I141").Value = 855.875
$ws.Range("J141").Value = 1500
$ws.Range("K141").Value = 2567.625
$ws.Range("L141").Value = 4500
$ws.Range("M141").Value = 2612.375
$ws.Range("N141").Value = -14860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 993.375
$ws.Range("I97").Value = 917.7273
$ws.Range("J97").Value = 1159.8
$ws.Range("K97").Value = 917.7273
$ws.Range("L97").Value = 1159.8
$ws.Range("M97").Value = -421.7273
$ws.Range("N97").Value = -2151.8

$ws.Range("H102").Value = 1577
$ws.Range("I102").Value = 1472.25
$ws.Range("J102").Value = 1716.6666
$ws.Range("K102").Value = 1472.25
$ws.Range("L102").Value = 1716.6666
$ws.Range("M102").Value = 149.75
$ws.Range("N102").Value = -4960.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1354.6666
$ws.Range("I100").Value = 1185
$ws.Range("J100").Value = 2033.3334
$ws.Range("K100").Value = 1185
$ws.Range("L100").Value = 2033.3334
$ws.Range("M100").Value = -644
$ws.Range("N100").Value = -3115.3334

$ws.Range("H122").Value = 5549.65
$ws.Range("I122").Value = 5587.5293
$ws.Range("J122").Value = 5335
$ws.Range("K122").Value = 16762.5879
$ws.Range("L122").Value = 16005
$ws.Range("M122").Value = -14312.5879
$ws.Range("N122").Value = -20905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1301.9131
$ws.Range("I107").Value = 1070.5714
$ws.Range("J107").Value = 1661.7778
$ws.Range("K107").Value = 3211.7142
$ws.Range("L107").Value = 4985.3334
$ws.Range("M107").Value = -1291.7142
$ws.Range("N107").Value = -8825.3334

$ws.Range("H113").Value = 481.94446
$ws.Range("I113").Value = 449.33334
$ws.Range("J113").Value = 514.55554
$ws.Range("K113").Value = 1348.00002
$ws.Range("L113").Value = 1543.66662
$ws.Range("M113").Value = 821.9999800000001
$ws.Range("N113").Value = -5883.66662
